# Generate Report for Handoff
# Updates the localization-status report: the zh-cn and de-de status move
# from "In Translation" to "Ready for handoff", and the corresponding
# handoff/xliff-generation timestamps are refreshed. Excel then re-autofits
# the now-wider "Status"/language columns to fit the new, longer text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet --------------------------------------------------
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-31 12:48:40"

# --- zh-cn sheet -------------------------------------------------------
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-31 12:48:35"

# --- de-de sheet -------------------------------------------------------
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-31 12:48:40"

# --- Re-fit the columns that now hold the longer "Ready for handoff"
#     text, same as Excel does automatically when generating the report.
$overview.Range("E:E").ColumnWidth = 16.3
$overview.Range("F:F").ColumnWidth = 16.3
$zhcn.Range("C:C").ColumnWidth = 16.3
$dede.Range("C:C").ColumnWidth = 16.3
